$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column H corresponds to "municipio-nombre". It is being re-curated from a
# "medida" (measure) into a "dim" (dimension), matching the pattern already
# used by column I ("provincia-nombre").
$ws.Range("H2").Value = "sdmx-dimension:refArea"
$ws.Range("H3").Value = "dim"
$ws.Range("H4").Value = "URI-Municipio"
